$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos list (GitHub Actions scheduled update): new price
# (column D) and volume-change (column E) readings for the existing coins,
# plus three coins (Stellar / Fetch.AI / ApeXProtocol / Monero / FirstDigitalUSD)
# whose relative ranking shuffled, so rows 45-49 get their Coin/Link/Price/
# Volume columns rewritten in the new order (the rank numbers in column A
# are unaffected).

# Force text formatting on cells whose new values look numeric,
# so Excel stores them as text (matching the source data) instead of
# silently converting to a Double (which would lose formatting like
# trailing zeros, e.g. '1.00', or multi-dot strings like '3.365.58').
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated values
$ws.Range("D2").Value = '64.925.72'
$ws.Range("E2").Value = '  +0.94%  '
$ws.Range("D3").Value = '3.365.58'
$ws.Range("E3").Value = '  +0.89%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '554.41'
$ws.Range("E5").Value = '  +0.51%  '
$ws.Range("D6").Value = '174.01'
$ws.Range("E6").Value = '  -0.54%  '
$ws.Range("D7").Value = '0.629'
$ws.Range("E7").Value = '  +2.03%  '
$ws.Range("D8").Value = '3.357.43'
$ws.Range("E8").Value = '  +0.92%  '
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("D10").Value = '0.174'
$ws.Range("E10").Value = '  +6.81%  '
$ws.Range("D11").Value = '0.635'
$ws.Range("E11").Value = '  +1.58%  '
$ws.Range("D12").Value = '53.52'
$ws.Range("E12").Value = '  -1.78%  '
$ws.Range("E13").Value = '  +3.42%  '
$ws.Range("D14").Value = '9.14'
$ws.Range("E14").Value = '  +1.16%  '
$ws.Range("D15").Value = '3.908.37'
$ws.Range("E15").Value = '  +0.98%  '
$ws.Range("E16").Value = '  +2.31%  '
$ws.Range("D17").Value = '18.22'
$ws.Range("E17").Value = '  -0.35%  '
$ws.Range("D18").Value = '3.358.75'
$ws.Range("E18").Value = '  +0.51%  '
$ws.Range("D19").Value = '64.984.92'
$ws.Range("E19").Value = '  +1.21%  '
$ws.Range("D20").Value = '11.84'
$ws.Range("E20").Value = '  +0.97%  '
$ws.Range("D21").Value = '0.995'
$ws.Range("E21").Value = '  +1.93%  '
$ws.Range("D22").Value = '451.26'
$ws.Range("E22").Value = '  +4.35%  '
$ws.Range("D23").Value = '4.92'
$ws.Range("E23").Value = '  -3.31%  '
$ws.Range("D24").Value = '4.07'
$ws.Range("E24").Value = '  +0.40%  '
$ws.Range("D25").Value = '86.96'
$ws.Range("E25").Value = '  +3.30%  '
$ws.Range("D26").Value = '13.70'
$ws.Range("E26").Value = '  +2.32%  '
$ws.Range("D27").Value = '10.75'
$ws.Range("E27").Value = '  +0.25%  '
$ws.Range("D28").Value = '2.87'
$ws.Range("E28").Value = '  +1.54%  '
$ws.Range("D29").Value = '8.63'
$ws.Range("E29").Value = '  -0.99%  '
$ws.Range("D30").Value = '30.99'
$ws.Range("E30").Value = '  +4.43%  '
$ws.Range("D31").Value = '6.55'
$ws.Range("E31").Value = '  -1.38%  '
$ws.Range("D32").Value = '62.92'
$ws.Range("E32").Value = '  +8.16%  '
$ws.Range("D33").Value = '11.43'
$ws.Range("E33").Value = '  -0.34%  '
$ws.Range("D34").Value = '576.10'
$ws.Range("E34").Value = '  -0.75%  '
$ws.Range("D35").Value = '0.107'
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("E37").Value = '  +3.87%  '
$ws.Range("D38").Value = '0.140'
$ws.Range("E38").Value = '  -0.29%  '
$ws.Range("D39").Value = '35.60'
$ws.Range("E39").Value = '  +0.03%  '
$ws.Range("D40").Value = '0.370'
$ws.Range("E40").Value = '  +1.08%  '
$ws.Range("D41").Value = '0.0₃0739'
$ws.Range("E41").Value = '  -1.52%  '
$ws.Range("D42").Value = '3.086.46'
$ws.Range("E42").Value = '  -0.54%  '
$ws.Range("D43").Value = '0.0417'
$ws.Range("E43").Value = '  +2.22%  '
$ws.Range("D44").Value = '2.76'
$ws.Range("E44").Value = '  -1.14%  '
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").Value = '0.134'
$ws.Range("E45").Value = '  +3.08%  '
$ws.Range("B46").Value = 'Fetch.AI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D46").Value = '2.45'
$ws.Range("E46").Value = '  -0.45%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").Value = '3.17'
$ws.Range("E47").Value = '  -0.80%  '
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").Value = '142.62'
$ws.Range("E48").Value = '  +5.67%  '
$ws.Range("B49").Value = 'FirstDigitalUSD'
$ws.Range("C49").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D49").Value = '1.00'
$ws.Range("E49").Value = '  +0.13%  '
$ws.Range("E50").Value = '  -2.25%  '
$ws.Range("D51").Value = '8.27'
$ws.Range("E51").Value = '  +0.34%  '

# Reset style back to Normal on the forced-text cells so no stray
# number-format style lingers on them
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
